# "Últimas simulações para o parâmetro de impacto"
# Update the simulation inputs (B2: lat, C2: semiEixo, D2: auxSemiEixo) on
# Sheet1 with the latest values; the dependent formulas in E2 (auxDegree)
# and K2 (angInc) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 79.57
$ws.Range("C2").Value = 1.59
$ws.Range("D2").Value = 0.038

# Move the active selection to B3, matching the latest saved view.
[void]$ws.Range("B3").Select()
